# EIA Table 2.7.A monthly update: October 2016 -> November 2016 data refresh.
# Adds the new "November" monthly row, shifts the trailing summary rows down
# by one, refreshes the Annual Totals / Year to Date figures, and updates the
# title/footer text that references the reporting month.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- 1. Update the title/subtitle text that references the reporting month ---
$ws.Range("A2").Value = "by Sector, 2006-November 2016 (Billion Btus)"

# --- 2. Insert the new "November" data row right after the October row (53) ---
$ws.Rows("53:53").Insert()

# Match the formatting of the row above (the October data row) exactly, since
# a bare row Insert() does not carry the border formatting through.
$ws.Range("A52:F52").Copy()
$ws.Range("A53:F53").PasteSpecial(-4122)
$excel.CutCopyMode = 0

$ws.Range("A53").Value = "November"
$ws.Range("B53").Value = 29901
$ws.Range("C53").Value = 2839
$ws.Range("D53").Value = 11531
$ws.Range("E53").Value = 30
$ws.Range("F53").Value = 15500

# --- 3. Refresh the "Annual Totals" block (now rows 55-57) ---
$ws.Range("B55").Value = 392989
$ws.Range("C55").Value = 41568
$ws.Range("D55").Value = 158971
$ws.Range("E55").Value = 910
$ws.Range("F55").Value = 191539

$ws.Range("B56").Value = 371750
$ws.Range("C56").Value = 40169
$ws.Range("D56").Value = 156449
$ws.Range("E56").Value = 469
$ws.Range("F56").Value = 174664

$ws.Range("B57").Value = 347179
$ws.Range("C57").Value = 38758
$ws.Range("D57").Value = 139632
$ws.Range("E57").Value = 610
$ws.Range("F57").Value = 168179

# --- 4. Update the "Rolling 12 Months Ending in ..." label (now row 58) ---
$ws.Range("A58").Value = "Rolling 12 Months Ending in November"

# --- 5. Refresh the "Year to Date" block (now rows 59-60) ---
$ws.Range("B59").Value = 410046
$ws.Range("C59").Value = 44244
$ws.Range("D59").Value = 171990
$ws.Range("E59").Value = 520
$ws.Range("F59").Value = 193292

$ws.Range("B60").Value = 382080
$ws.Range("C60").Value = 42508
$ws.Range("D60").Value = 154570
$ws.Range("E60").Value = 646
$ws.Range("F60").Value = 184356

Write-Host "edit complete"
